$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.420.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.642.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.536"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.32%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.09"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("E10").Value = "  -2.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0890"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.875.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.649.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("E14").Value = "  -3.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.556"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.394.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0719"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.66%  "
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  -3.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("E26").Value = "  +1.91%  "
$ws.Range("E27").Value = "  -3.39%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.34%  "
$ws.Range("E30").Value = "  -4.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0486"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.17%  "
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.412.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.45%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  -1.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.881"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.14%  "
$ws.Range("E39").Value = "  -3.40%  "
$ws.Range("E40").Value = "  +0.99%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  -1.55%  "
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.12%  "
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.786.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("E48").Value = "  -3.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.11%  "
$ws.Range("E50").Value = "  -3.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0987"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.66%  "
